$wb = $excel.ActiveWorkbook
$members = $excel | Get-Member
$members | ForEach-Object { Write-Output $_.Name }
